$wb = $excel.ActiveWorkbook

# --- GUID / identifier substitutions -------------------------------------
$oldGuid1 = "c43b3f1f-b159-4eef-8940-5e41d1f38fc6"
$newGuid1 = "cb8b6387-ef20-4e9d-963c-e493c86ea27a"
$oldGuid2 = "e2944927-23b0-4647-9543-ab3703d28b1e"
$newGuid2 = "ffff9ba7ca40-387a-4f0e-8e98-85c9db8fb613"

$oldHash  = "f96392784b74c1d958608579280229792a221736"
$newHash  = "379c70dac170965e32cac0c420ccc280e6549c9b"

# --- Timestamp substitutions ----------------------------------------------
$tsOverview_old = "2016-08-26 09:07:11"
$tsOverview_new = "2016-08-26 09:08:30"

$tsZhHandoff_old  = "2016-08-26 09:06:59"
$tsZhHandoff_new  = "2016-08-26 09:08:25"
$tsZhHandback_old = "2016-08-26 09:07:30"
$tsZhHandback_new = "2016-08-26 09:08:42"

$tsDeHandoff_old  = "2016-08-26 09:07:11"
$tsDeHandoff_new  = "2016-08-26 09:08:30"
$tsDeHandback_old = "2016-08-26 09:07:37"
$tsDeHandback_new = "2016-08-26 09:08:49"

# --- Sheet: Overview --------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid1.md"
$wsOverview.Range("G2").Value = $tsOverview_new

$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newGuid2.md"
$wsOverview.Range("G3").Value = $tsOverview_new

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.TextToDisplay -eq "e2e\$oldGuid1.md") {
        $hl.TextToDisplay = "e2e\$newGuid1.md"
    } elseif ($hl.TextToDisplay -eq "e2e\$oldGuid2.md") {
        $hl.TextToDisplay = "e2e\$newGuid2.md"
    }
}

# --- Sheet: zh-cn ------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid1.md"
$wsZh.Range("G2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = $tsZhHandoff_new
$wsZh.Range("I2").Value = "$newGuid1.md"
$wsZh.Range("J2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZh.Range("K2").Value = $tsZhHandback_new

$wsZh.Range("A3").Value = "$newGuid2.md"
$wsZh.Range("G3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZh.Range("H3").Value = $tsZhHandoff_new
$wsZh.Range("I3").Value = "$newGuid2.md"
$wsZh.Range("J3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZh.Range("K3").Value = $tsZhHandback_new

foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.TextToDisplay -eq "$oldGuid1.md") {
        $hl.TextToDisplay = "$newGuid1.md"
    } elseif ($hl.TextToDisplay -eq "$oldGuid2.md") {
        $hl.TextToDisplay = "$newGuid2.md"
    }
}

# --- Sheet: de-de ------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid1.md"
$wsDe.Range("G2").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = $tsDeHandoff_new
$wsDe.Range("I2").Value = "$newGuid1.md"
$wsDe.Range("J2").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDe.Range("K2").Value = $tsDeHandback_new

$wsDe.Range("A3").Value = "$newGuid2.md"
$wsDe.Range("G3").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDe.Range("H3").Value = $tsDeHandoff_new
$wsDe.Range("I3").Value = "$newGuid2.md"
$wsDe.Range("J3").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDe.Range("K3").Value = $tsDeHandback_new

foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.TextToDisplay -eq "$oldGuid1.md") {
        $hl.TextToDisplay = "$newGuid1.md"
    } elseif ($hl.TextToDisplay -eq "$oldGuid2.md") {
        $hl.TextToDisplay = "$newGuid2.md"
    }
}

Write-Output "done"
